$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update summary figures -------------------------------------------------
# VALOR MORA (total)
$ws.Range("E11").Value2 = 215865
# Cant. Trabajadores
$ws.Range("C13").Value2 = 4
# Cant. Periodos
$ws.Range("F13").Value2 = 3

# --- Drop three of the old employee detail rows -----------------------------
# Removing rows 16-18 shifts the remaining detail row (old row 19-22, with the
# distinct "last row" formatting that lived on row 22) up to rows 16-19, and
# shifts the signature block up from rows 27/28 to rows 24/25.
$ws.Rows("16:18").Delete()

# --- Populate the (now four) employee detail rows with the new dataset -----
# Row 16: MARGENIA DE JESUS BURGOS PEREZ
$ws.Range("B16").Value2 = "CC"
$ws.Range("C16").Value2 = "26162681"
$ws.Range("D16").Value2 = "MARGENIA DE JESUS BURGOS PEREZ"
$ws.Range("E16").Value2 = "1911"
$ws.Range("F16").Value2 = 31249
$ws.Range("G16").Value2 = 781242

# Row 17: ALEXANDER REALES RAMIREZ
$ws.Range("B17").Value2 = "CC"
$ws.Range("C17").Value2 = "11165994"
$ws.Range("D17").Value2 = "ALEXANDER REALES RAMIREZ"
$ws.Range("E17").Value2 = "2003"
$ws.Range("F17").Value2 = 59348
$ws.Range("G17").Value2 = 1483712

# Row 18: JORGE LUIS QUINTANA MARTINEZ
$ws.Range("B18").Value2 = "CC"
$ws.Range("C18").Value2 = "3805709"
$ws.Range("D18").Value2 = "JORGE LUIS QUINTANA MARTINEZ"
$ws.Range("E18").Value2 = "2505"
$ws.Range("F18").Value2 = 68328
$ws.Range("G18").Value2 = 1708200

# Row 19: WALTER DE LA CRUZ ASENCIO CHAMORRO
$ws.Range("B19").Value2 = "CC"
$ws.Range("C19").Value2 = "9314008"
$ws.Range("D19").Value2 = "WALTER DE LA CRUZ ASENCIO CHAMORRO"
$ws.Range("E19").Value2 = "2505"
$ws.Range("F19").Value2 = 56940
$ws.Range("G19").Value2 = 1423500
